# Actualización automática hashcode mar jul 16 02:08:13 CEST 2019
# Updates the "hashcode" values (column B) for specific "code" rows (column A)
# on the "hashcode.csv" worksheet, matching the target commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @(
    @{ Row = 11;  Old = "1352d9b99bf06626ff80952eda02d7d2"; New = "1f682c4baf00039722b9d3b2a8f6431f" },
    @{ Row = 34;  Old = "c61e0c5fa0c3d3aeb7f195c62229f494"; New = "9b5fa738b68a8c46f512c3e8ae609d3b" },
    @{ Row = 44;  Old = "a2cfcbfef9b7b4aed5ed06cdf76e820f"; New = "775da89266fde57dfe7ca7c89abf5d91" },
    @{ Row = 74;  Old = "9555bf74da8a390313ded720eb47dce7"; New = "8a74666dc4ebb183229cedc771aa374f" },
    @{ Row = 89;  Old = "160ee88f449d69ffbf488ebe9d2dcc44"; New = "e5a9c26e094a5557ae9c4aa83e416d55" },
    @{ Row = 99;  Old = "ec5bd2a050b8a245967e920be6cdaaa2"; New = "0c473cacc596f7b80f753639d0d0ca9c" },
    @{ Row = 110; Old = "4050bd447a74401c61ea746f9711d4fc"; New = "8c9098805d070995ea6995c660cc73a1" },
    @{ Row = 121; Old = "27c1bb70cb640d5ca20a759347c927c8"; New = "81667d4f5140992663fc6287a415e11f" },
    @{ Row = 154; Old = "e9828e955ed4896624069e2230da5da2"; New = "0164192226833e8b2508d9634b0ba903" },
    @{ Row = 160; Old = "f3de5288eeaf606f566c40f38f1f948a"; New = "adf3c1215f1ec05392a34e4fcab6d818" },
    @{ Row = 161; Old = "9bb4c7968671c6ffbee5b3db18131f17"; New = "1e5c3f3bf56fea72588394470e1cc359" },
    @{ Row = 162; Old = "28b7081ddd8b2bf574091a34d8703cef"; New = "537a5222143850acb0b8e7c2a56d1a6f" },
    @{ Row = 168; Old = "36c8cd53ba8a46717318adc0a51706b1"; New = "bc95cae257a5ff8399d8aa38ac0096e0" },
    @{ Row = 180; Old = "4452182d4a3e39871668d09fdb6c1e5b"; New = "8e3e66726412138b9c21d57bc4009d98" },
    @{ Row = 191; Old = "c73e5ad0a567948972aa3db3a087d497"; New = "aec159b771e496e8cb54e48f8a239e8e" },
    @{ Row = 213; Old = "e11742ebab986b101aaf472dd8371e81"; New = "f1a3da6a4991d211f4d0e18b9486ed7a" },
    @{ Row = 278; Old = "4f4e6e1d7f91885a3a4f184b8ac396e3"; New = "9283cf6e227051ed64790cd8214746ac" },
    @{ Row = 293; Old = "21201fdc44ce87e98d9209da669acf6b"; New = "a7d0b31354aa502f18e0103883abbc31" },
    @{ Row = 335; Old = "ecbe729ac86df7acbe5e7934836f2f14"; New = "fa67257d9e82773e7b9d6f5b58515c14" },
    @{ Row = 345; Old = "183913fecc02620ae6913e0667b17656"; New = "3d3502f758d76be92c0f4e2ea3201dd1" },
    @{ Row = 461; Old = "b11b80ec3b93464d6b97a5f9c1948435"; New = "060072cb4a449d58d07838c00b609f70" },
    @{ Row = 480; Old = "f23b3dca7b162c63f81a3379142179f4"; New = "1fd9ef0f8869fc52d6c81138b24ec41c" },
    @{ Row = 506; Old = "51d94fbb108c060af0774f3dfc25fd2e"; New = "aa1791820592e49d2dde3aff5748084a" },
    @{ Row = 514; Old = "1ff4dd27e25e4cecffa8c888a063c5c2"; New = "0163ad4ebad868ebcb1fb1d515410e6b" },
    @{ Row = 524; Old = "586802b4d9ba45de50d961c63708f3c0"; New = "b8463e643f40c14c051b7aa3e19cc647" },
    @{ Row = 534; Old = "76da3783aa2a61aa6867b6ba825b3179"; New = "b4d216af1c0225064ccc574065e16246" },
    @{ Row = 547; Old = "12134a6651c6de21c72dc6c1e1dae89a"; New = "61c4f18193adac7d146bc75c0f680430" },
    @{ Row = 553; Old = "58d85ba2051dd71507a5e4255d2e5b94"; New = "8317bc5e1079993b6d686cc7d773b4ef" },
    @{ Row = 572; Old = "f1eff8d1240251c266d684e4cbc1fca7"; New = "5ed55f8b2ae0bd9cea467720286f267b" },
    @{ Row = 584; Old = "a576e1b2662d1a21d6c1d37626fd4452"; New = "e375d004872e7eac94fce210d9414135" },
    @{ Row = 666; Old = "6a504f8d367e29df8fe91b6e061f2350"; New = "d0198b482e7ad0701fea272aba6657a8" },
    @{ Row = 729; Old = "27ed38bf1fbffac7273df8279ccba7ca"; New = "b4db0bd5cfe9f51ea71702c7935a8b82" },
    @{ Row = 768; Old = "8a866f38cea4d509d812189b47eef642"; New = "856d009b685edcaa25e7aebd1e4cb92c" },
    @{ Row = 811; Old = "dbd952bba9bedbb15ced3d14a76bc9b0"; New = "5f1e48ea2ee37ac4a0cd6534daf28e1d" },
    @{ Row = 815; Old = "bd5b9380588c9dc7c9ba8123dc3cab76"; New = "deeeabb02d47e448e34e5d3bbaeb8dad" },
    @{ Row = 816; Old = "1951623ae9020a139ec3467817acc2ab"; New = "831b12f239db1883cfb6a62cd480eabe" },
    @{ Row = 825; Old = "76fb08e3968f1341beee8c4d704ab1a6"; New = "e0b748b7abab51601ff88878e1646e1d" },
    @{ Row = 827; Old = "fe391b223dd9b3e7fc6a5f6ebd9890a3"; New = "e72e4ad52475855fd285dd2b5bbecbd4" },
    @{ Row = 874; Old = "d878f735a89572d2273c1e98708e28dd"; New = "c9c849f03081bb7a17b5eba5feebb7ea" }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, 2)
    if ($cell.Value2 -ne $u.Old) {
        Write-Host ("WARNING: row {0} column B expected '{1}' but found '{2}'" -f $u.Row, $u.Old, $cell.Value2)
    }
    $cell.Value2 = $u.New
}

Write-Host ("Updated {0} hashcode cells." -f $updates.Count)
